$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1) and "全部类型" sheet (sheet4) with new F-column values
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 652
    $ws.Range("F3").Value = 3885
    $ws.Range("F5").Value = 729
}
